# Config-file-format manual: add EyeLink/Optotrak recording fields
# (trial_timeout_msg, trial_kb_resp, text_color) and update several
# field descriptions. This shifts the old rows 8-29 down to 11-32 and
# fills rows 8-32 with the new documentation table content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for three new rows (trial_timeout_msg, trial_kb_resp, text_color)
# above the old "background_color" row, pushing everything below down.
$ws.Rows("8:10").Insert()

$data = @(
    @(8,  "trial_timeout_msg",      "bool",   "N/A",           "0=no timeout message, 1 = display timeout message if trial ends with no response"),
    @(9,  "trial_kb_resp",          "bool",   "N/A",           "0 = no keyboard response, 1 = prompt for keyboard response at trial end"),
    @(10, "trial_feedback",         "string", "N/A",           "image file name, or ""msg"" which will display the message specified at the beginning of the script"),
    @(11, "background_color",       "string", "RGB (0-1)",     "examples:  ""0 0 0"" = black, ""1 1 1"" = white, ""0.5 0.5 0.5"" = grey, ""1 0 0"" = red, etc. [default = white]"),
    @(12, "text_color",             "string", "RGB (0-1)",     "examples:  ""0 0 0"" = black, ""1 1 1"" = white, ""0.5 0.5 0.5"" = grey, ""1 0 0"" = red, etc. [default = black]"),
    @(13, "stim_img_name",          "string", "N/A",           "name of image file (typically in the ""Images/"" subfolder)"),
    @(14, "stim_onset",             "float",  "seconds",       "stimulus onset time (from trial start time, EXCLUDING fixation pauses)"),
    @(15, "stim_duration",          "float",  "seconds",       "duration that stimulus image is displayed onscreen"),
    @(16, "stim_cent_x",            "float",  "percent (0-1)", "stimulus center (percentage of screen width)"),
    @(17, "stim_cent_y",            "float",  "percent (0-1)", "stimulus center (percentage of screen height)"),
    @(18, "stim_size_x",            "float",  "cm",            "stimulus actual size on screen"),
    @(19, "stim_size_y",            "float",  "cm",            "stimulus actual size on screen"),
    @(20, "stim_rotation",          "float",  "degrees",       "stimulus rotation"),
    @(21, "stim_is_touchable",      "int",    "N/A",           "0 = not touchable, 1 = touch only, 2 = gaze only,  3 = touch or gaze, 4 = touch and gaze"),
    @(22, "stim_is_target",         "int",    "N/A",           "0 = not target, 1 = touch target, 2 =  gaze target"),
    @(23, "subj_fixation_type",     "int",    "N/A",           "0 = not a fixation object, 1 = touch fixation, 2 = gaze fixation"),
    @(24, "subj_fixation_onset",    "float",  "seconds",       "fixation onset time (from trial start time, INCLUDING fixation pauses)"),
    @(25, "subj_fixation_duration", "float",  "seconds",       "required consecutive fixation time by participant"),
    @(26, "mask_onset",             "float",  "seconds",       "stimulus mask onset time (from trial start time, EXCLUDING fixation pauses)"),
    @(27, "mask_duration",          "float",  "seconds",       "duration that stimulus mask is displayed onscreen"),
    @(28, "mask_size",              "float",  "cm",            "stimulus mask dot size"),
    @(29, "mask_color",             "float",  "N/A",           "stimulus mask dot color (0 = black, 1=white, others are possible as an RGB vector [r g b] floats 0-1)"),
    @(30, "mask_rotation",          "float",  "degrees",       "stimulus mask rotation"),
    @(31, "mask_fit",               "int",    "N/A",           "0 = mask is bounding rectangle, 1 = mask is fitted along detected shape borders of stim image"),
    @(32, "mask_margin",            "float",  "cm",            "margin of mask points around stim image")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# subj_fixation_type row ended up very slightly shorter than the default
$ws.Rows("23:23").RowHeight = 14.25

$ws.Range("D22").Select() | Out-Null
